$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J header: SCE loop moved from (Q)uarterly to (M)onthly estimate.
$ws.Range("J1").Value = 'SE: $\hat\lambda_{SCE}$(M)'

# Bounded loop estimates for column I (lambda SPF) are now stable at 1
# instead of drifting (1.8, 1.74, 1.77) for the unbounded rows.
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
